# Server now generates translations; rename the "name" header/key column
# to "_name" so it isn't confused with a translatable string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "_name"
